$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.759.90"
$ws.Range("E2").Value = "  -1.82%  "
$ws.Range("D3").Value = "1.869.15"
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "300.16"
$ws.Range("E5").Value = "  -2.37%  "
$ws.Range("E6").Value = "  -0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5320"
$ws.Range("E7").Value = "  +1.12%  "
$ws.Range("E8").Value = "  -2.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07153"
$ws.Range("E9").Value = "  -1.81%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "21.44"
$ws.Range("E10").Value = "  -2.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8857"
$ws.Range("E11").Value = "  -1.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08154"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("D13").Value = "1.877.86"
$ws.Range("E13").Value = "  +28.73%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "92.37"
$ws.Range("E14").Value = "  -3.74%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.286"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.81"
$ws.Range("E17").Value = "  +0.36%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008482"
$ws.Range("E18").Value = "  -1.98%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.000"
$ws.Range("E19").Value = "  -0.05%  "
$ws.Range("D20").Value = "26.795.04"
$ws.Range("E20").Value = "  -1.82%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.968"
$ws.Range("E21").Value = "  -2.90%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.62"
$ws.Range("E22").Value = "  -1.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.370"
$ws.Range("E23").Value = "  -2.26%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.276"
$ws.Range("E24").Value = "  -1.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "145.69"
$ws.Range("E25").Value = "  -2.75%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.731"
$ws.Range("E26").Value = "  -0.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.00"
$ws.Range("E27").Value = "  -1.39%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "113.51"
$ws.Range("E28").Value = "  -2.94%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "4.692"
$ws.Range("E29").Value = "  -3.21%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.630"
$ws.Range("E30").Value = "  -4.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.09102"
$ws.Range("E31").Value = "  -1.55%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8074"
$ws.Range("E32").Value = "  -2.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05011"
$ws.Range("E33").Value = "  -1.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.172"
$ws.Range("E34").Value = "  -4.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.944"
$ws.Range("E35").Value = "  -1.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6113"
$ws.Range("E36").Value = "  +5.49%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.641"
$ws.Range("E37").Value = "  -2.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.178"
$ws.Range("E38").Value = "  -5.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01941"
$ws.Range("E39").Value = "  -3.16%  "
$ws.Range("E40").Value = "  -1.39%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.5253"
$ws.Range("E41").Value = "  +6.62%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.468"
$ws.Range("E42").Value = "  -1.91%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.708"
$ws.Range("E43").Value = "  -5.66%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "115.14"
$ws.Range("E44").Value = "  -1.38%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.1489"
$ws.Range("E45").Value = "  -2.26%  "
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.638"
$ws.Range("E47").Value = "  -0.60%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "9.902"
$ws.Range("E48").Value = "  -2.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "37.27"
$ws.Range("E49").Value = "  -4.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06060"
$ws.Range("E50").Value = "  -1.91%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "62.04"
$ws.Range("E51").Value = "  -3.93%  "
